# Apply the "instance_name setting" edit to household_member.xlsx
#
# Summary of the change:
#  - survey!D6 (the stop_survey-ish note prompt) drops the trailing
#    " for {{instance_name}}" from its label text.
#  - settings sheet gains a new row: instance_name = member_name,
#    documenting the new setting that replaces the inline {{instance_name}}.
#  - survey row 6 grows taller (the note label got shorter, but the row was
#    re-wrapped / resized by hand afterwards).
#  - Active sheet/selection bookkeeping: "settings" becomes the active tab
#    (with a new selection at B18), "survey" is no longer the tab-selected
#    sheet (selection moves to D7).

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- survey sheet: shorten the note label, resize row 6, move selection ---
$survey.Range("D6").Value = "{{member_name}} age is {{evaluate calculates.ageIsOddOrEven}} in {{setting 'table_id'}}"
$survey.Rows.Item(6).RowHeight = 62.5
$survey.Range("D7").Select() | Out-Null

# --- settings sheet: add the new instance_name setting row ---
$settings.Range("A6").Value = "instance_name"
$settings.Range("B6").Value = "member_name"

# --- activate the settings sheet last so it becomes the active tab, with
#     its own new selection ---
$settings.Activate()
$settings.Range("B18").Select() | Out-Null

$wb.Save()
